# Update the "丽水-漫展信息" workbook: both the "展览" sheet and the
# "全部类型" sheet carry the same exhibition table and both received the
# same edit upstream (new row for the 龙泉 event inserted as row 3, pushing
# the 丽水·LPJ row down to row 4, plus updated "想去人数" counters).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # "想去人数" (want-to-go count) bump on the first event row.
    $ws.Range("F2").Value = 252

    # Insert a brand-new row 3 for the 龙泉 event; this pushes the existing
    # row 3 (丽水·LPJ ...) down to row 4 and grows the sheet dimension to
    # A1:J4 automatically.
    $ws.Rows.Item(3).Insert()

    # --- new row 3: 龙泉·崩X铁X原ONLY -------------------------------------
    $ws.Range("A3").Value = 2
    $ws.Range("A3").Font.Bold = $true
    $ws.Range("A3").HorizontalAlignment = -4108
    $ws.Range("A3").VerticalAlignment = -4160
    $ws.Range("A3").Borders.LineStyle = 1

    $ws.Range("B3").Value = "'2024.02.07"
    $ws.Range("C3").Value = "龙泉·崩X铁X原ONLY"
    $ws.Range("D3").Value = "金沙路26-1号 龙泉金沙温泉酒店"
    $ws.Range("E3").Value = "'2024.02.07 10:30-02.07 16:30"
    $ws.Range("F3").Value = 0
    $ws.Range("G3").Value = "'50"
    $ws.Range("H3").Value = $false
    $ws.Range("I3").Value = "https://show.bilibili.com/platform/detail.html?id=80714&msource=Msearch_colligation"
    $ws.Range("J3").Value = "//i2.hdslb.com/bfs/openplatform/202401/rTvQio211704877379770.jpeg"

    # --- row 4 (formerly row 3, shifted down): renumber + refresh count ---
    $ws.Range("A4").Value = 3
    $ws.Range("F4").Value = 249
}
